# Language workbook update - stage 1 body changes
#
# Removes the "Tiny Bacillus" (bodyBacillusTiny) and "Tiny Coccus"
# (bodyCoccusTiny) body-type rows from the sheet (these keys/values are
# dropped from the shared strings table automatically once no cell
# references them), shifting the remaining rows up to close the gaps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Tiny Coccus" row (currently row 17) first, then the
# "Tiny Bacillus" row (currently row 15), so row numbers for the first
# deletion aren't invalidated by the second.
$ws.Rows(17).Delete()
$ws.Rows(15).Delete()

# Update the view/selection to match the post-edit cursor position.
$ws.Range("B16").Select()
